# Excel COM-interop script recreating the authored edit:
#  1. Rename "Feuil1" -> "data"
#  2. Move the helper "allocation / produit structuré" table that lived in
#     columns I:T of the data sheet into a brand-new "temp" worksheet
#     (as columns A:L), and remove it from "data".
#  3. Turn the data range A1:G3017 on "data" into an AutoFilter table and
#     register the hidden _xlnm._FilterDatabase name (scoped to "data").
#  4. Fix 4 cells (G1003:G1006) that held the literal Bloomberg placeholder
#     text "#N/A N/A" so they hold the real numeric value -0.451 instead.
#  5. Restore sensible view state (active sheet/selection) on both sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "data"

# --- 2. Create the "temp" sheet right after "data" and move I1:T21 there ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$ws2.Name = "temp"

$helper = $ws.Range("I1:T21")
$helper.Copy($ws2.Range("A1"))

# Remove the helper table from "data": drop the now-empty K:T columns
# entirely, and clear out I:J (which stay part of the sheet but carry no
# more data after the move).
$ws.Range("K1:T3017").EntireColumn.Delete()
$ws.Range("I1:J3017").Delete()

# --- 4. Replace the "#N/A N/A" placeholders with the real figure ---
$ws.Range("G1003:G1006").Value2 = -0.451

# --- 3. AutoFilter + hidden _FilterDatabase name on "data" ---
$ws.Range("A1:G3017").AutoFilter()
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=data!`$A`$1:`$G`$3017")
$fdb.Visible = $false

# --- 5. View state: gridlines off everywhere, selections, active sheet ---
$excel.ActiveWindow.DisplayGridlines = $false
$ws2.Range("J22").Select()
$ws.Activate()
$ws.Range("G1002").Select()
